$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E3"  = 16.432
    "A12" = -21.626
    "E14" = 16.994
    "E26" = 16.41400000000001
    "E31" = 16.376
    "A32" = -21.744
    "E35" = 16.63
    "A36" = -20.339
    "E37" = 16.737
    "A38" = -19.741
    "E45" = 16.851
    "A46" = -21.835
    "A54" = -21.869
    "A55" = -22.164
    "E57" = 16.643
    "A67" = -21.565
    "A69" = -21.721
    "A72" = -21.445
    "A91" = -21.652
    "A99" = -20.43
    "E100" = 16.749
    "E102" = 16.669
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
